# Update the "dSF" column (F) values for specific rows to reflect the
# repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    4  = -2
    6  = -7
    8  = -2
    10 = -4
    11 = -5
    13 = 4
    14 = -7
    16 = 3
    19 = -2
    23 = -1
    28 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
